$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$paidRows = @(25,26,27,28,29,31,33,34,35,36,37)
foreach ($r in $paidRows) {
    $ws.Cells.Item($r, 8).Value = "y"
}
$ws.Cells.Item(30, 8).Value = ""
